$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark in the table cell (Holiday"s").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Remove the stray "s" run at the very start of the document (before
#    "Student Progress Report") and add a "_GoBack" bookmark in its place.
$firstPara = $d.Paragraphs(1).Range
$firstPara.Find.Execute("s", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sRange = $d.Range($firstPara.Start, $firstPara.Start + 1)
$sRange.Text = ""
$d.Bookmarks.Add("_GoBack", $sRange)
